$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2025-12-29 13:57:30", "Admin", "quotation", "access_granted", "Opened quotation page"),
    @("2025-12-29 13:57:37", "Admin", "quotation", "access_granted", "Opened quotation page"),
    @("2025-12-29 13:57:39", "Admin", "quotation", "access_granted", "Opened quotation page"),
    @("2025-12-29 13:57:40", "Admin", "quotation", "access_granted", "Opened quotation page"),
    @("2025-12-29 13:57:40", "Admin", "quotation", "access_granted", "Opened quotation page"),
    @("2025-12-29 13:57:42", "Admin", "dashboard", "access_granted", "Opened dashboard page")
)

$startRow = 161
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
}
